$wb = $excel.ActiveWorkbook

# --- Original sheet becomes "0" -------------------------------------------
$ws0 = $wb.Worksheets.Item(1)

# Row-height cleanup on the original sheet: most wrapped rows collapse back
# to the single-line default (15.75pt); row 5 loses its explicit height
# entirely (falls back to the sheet default).
$rowsToReset = @(2,3,6,10,12,13,14,15,16,18,22,23,24,25,26,28)
foreach ($r in $rowsToReset) {
    $ws0.Rows.Item($r).RowHeight = 15.75
}
$ws0.Rows.Item(5).AutoFit() | Out-Null

# Selection on sheet "0" moves to B17, and it is no longer the active tab.
$ws0.Range("B17").Select() | Out-Null
$ws0.Name = "0"

# --- Add ten more sheets: "1".."10" ----------------------------------------
$prev = $ws0
$newSheets = @{}
for ($i = 1; $i -le 10; $i++) {
    $wsNew = $wb.Worksheets.Add($null, $prev)
    $wsNew.Name = "$i"
    $newSheets[$i] = $wsNew
    $prev = $wsNew
}

# A couple of the new sheets keep a non-default selection.
$newSheets[3].Range("E18").Select() | Out-Null

# Sheet "1" is the active / selected tab, selection on B1.
$newSheets[1].Activate() | Out-Null
$newSheets[1].Range("B1").Select() | Out-Null
